# Updated cryptos list with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) columns for the
# crypto rows with the latest scraped figures. All of these cells are
# stored as text (not numbers), so numeric-looking prices are written with
# a leading "'" quote-prefix to force Excel to keep them as text instead of
# silently re-typing them as numbers (which would also strip meaningful
# trailing zeros, e.g. "0.07680" -> 0.0768).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.192.23"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").Value = "1.905.99"

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'307.66"
$ws.Range("E5").Value = "  +0.95%  "

$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").Value = "'0.5242"
$ws.Range("E7").Value = "  +3.02%  "

$ws.Range("E8").Value = "  +3.02%  "

$ws.Range("D9").Value = "'0.07255"
$ws.Range("E9").Value = "  +1.02%  "

$ws.Range("E10").Value = "  +2.87%  "

$ws.Range("E11").Value = "  +0.56%  "

$ws.Range("D12").Value = "'0.07680"
$ws.Range("E12").Value = "  +2.50%  "

$ws.Range("D13").Value = "1.913.25"
$ws.Range("E13").Value = "  +2.22%  "

$ws.Range("D14").Value = "'95.35"
$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("D15").Value = "'5.274"
$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("D17").Value = "'0.000008590"
$ws.Range("E17").Value = "  +0.93%  "

$ws.Range("D18").Value = "'14.41"
$ws.Range("E18").Value = "  +1.80%  "

$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").Value = "27.261.42"

$ws.Range("D21").Value = "'5.062"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").Value = "2.151.16"
$ws.Range("E22").Value = "  +1.69%  "

$ws.Range("E23").Value = "  +2.45%  "

$ws.Range("D24").Value = "'6.428"
$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("D25").Value = "'2.304"
$ws.Range("E25").Value = "  +10.62%  "

$ws.Range("D26").Value = "'145.74"
$ws.Range("E26").Value = "  -1.53%  "

$ws.Range("E27").Value = "  -1.96%  "

$ws.Range("D28").Value = "'18.15"
$ws.Range("E28").Value = "  +1.54%  "

$ws.Range("D29").Value = "'114.74"
$ws.Range("E29").Value = "  +1.07%  "

$ws.Range("E30").Value = "  +5.36%  "

$ws.Range("D31").Value = "'4.797"
$ws.Range("E31").Value = "  +2.25%  "

$ws.Range("D32").Value = "'0.09215"
$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").Value = "'0.8087"
$ws.Range("E33").Value = "  +7.72%  "

$ws.Range("D34").Value = "'0.05049"
$ws.Range("E34").Value = "  -0.11%  "

$ws.Range("D35").Value = "'1.241"
$ws.Range("E35").Value = "  +7.56%  "

$ws.Range("D36").Value = "'2.997"
$ws.Range("E36").Value = "  +0.56%  "

$ws.Range("E37").Value = "  +3.02%  "

$ws.Range("E38").Value = "  +2.56%  "

$ws.Range("D39").Value = "'0.5677"
$ws.Range("E39").Value = "  +0.49%  "

$ws.Range("D40").Value = "'0.01982"
$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("D41").Value = "'1.074"
$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("D42").Value = "'8.980"
$ws.Range("E42").Value = "  +5.09%  "

$ws.Range("D43").Value = "'119.14"
$ws.Range("E43").Value = "  +3.19%  "

$ws.Range("D44").Value = "'6.617"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").Value = "'0.1512"
$ws.Range("E45").Value = "  +2.06%  "

$ws.Range("E46").Value = "  +0.99%  "

$ws.Range("D47").Value = "'10.18"
$ws.Range("E47").Value = "  +0.94%  "

$ws.Range("E48").Value = "  +0.18%  "

$ws.Range("E49").Value = "  +4.16%  "

$ws.Range("D50").Value = "'37.52"
$ws.Range("E50").Value = "  +1.45%  "

$ws.Range("D51").Value = "'63.67"
$ws.Range("E51").Value = "  +0.89%  "
